$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update comment text for speed sensor spacer 3mm (front hub, row 11) ---
# Was "to position WT_02004" -> now clarifies it became 2mm thick
$ws.Range("E11").Value = "to position WT_02004 became 2mm thick"
$ws.Range("F11").Formula = "=2*2"

# --- Update comment text for speed sensor spacer (front hub, row 9) ---
# Was "to position WT_02004" -> now clarifies the 1mm thickness
$ws.Range("E9").Value = "to position WT_02004 (1mm thick)"
$ws.Range("F9").Formula = "=2*1"

# --- New comment for rear hub Brake Bell (row 14) ---
$ws.Range("E14").Value = "from rear brake disc deal"

# --- New quantity formulas for rear hub speed sensor disc spacer rows ---
$ws.Range("F17").Formula = "=2*1"
$ws.Range("F19").Formula = "=2*2"

# --- Update the active selection to reflect where the editor left off ---
$ws.Range("E15").Select()
